$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in the title cell (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 19 de Marzo de 2020 a las 21:24"

# Rewrite the data rows (4-60) with the refreshed/reordered province data.
# This removes the duplicate "Cataluna*" row, moves Murcia/Melilla/La Palma to
# their new positions, and refreshes the case counts for the 21:24 snapshot.
$ws.Range("A4").Value2 = "Madrid"
$ws.Range("B4").Value2 = 6777
$ws.Range("C4").Value2 = 941
$ws.Range("D4").Value2 = 5338
$ws.Range("E4").Value2 = 498
$ws.Range("A5").Value2 = "Cataluña"
$ws.Range("B5").Value2 = 3270
$ws.Range("C5").Value2 = 3
$ws.Range("D5").Value2 = 3185
$ws.Range("E5").Value2 = 82
$ws.Range("A6").Value2 = "Araba/Alava"
$ws.Range("B6").Value2 = 621
$ws.Range("C6").Value2 = 18
$ws.Range("D6").Value2 = 583
$ws.Range("E6").Value2 = 35
$ws.Range("A7").Value2 = "Valencia/Valencia"
$ws.Range("B7").Value2 = 522
$ws.Range("C7").Value2 = 7
$ws.Range("D7").Value2 = 503
$ws.Range("E7").Value2 = 12
$ws.Range("A8").Value2 = "Navarra"
$ws.Range("B8").Value2 = 482
$ws.Range("C8").Value2 = 2
$ws.Range("D8").Value2 = 476
$ws.Range("E8").Value2 = 4
$ws.Range("A9").Value2 = "La Rioja"
$ws.Range("B9").Value2 = 468
$ws.Range("C9").Value2 = 2
$ws.Range("D9").Value2 = 461
$ws.Range("E9").Value2 = 5
$ws.Range("A10").Value2 = "Bizkaia/Vizcaya"
$ws.Range("B10").Value2 = 393
$ws.Range("C10").Value2 = 18
$ws.Range("D10").Value2 = 380
$ws.Range("E10").Value2 = 13
$ws.Range("A11").Value2 = "Malaga"
$ws.Range("B11").Value2 = 361
$ws.Range("C11").Value2 = 0
$ws.Range("D11").Value2 = 348
$ws.Range("E11").Value2 = 13
$ws.Range("A12").Value2 = "Alacant/Alicante"
$ws.Range("B12").Value2 = 338
$ws.Range("C12").Value2 = 11
$ws.Range("D12").Value2 = 316
$ws.Range("E12").Value2 = 11
$ws.Range("A13").Value2 = "Toledo"
$ws.Range("B13").Value2 = 293
$ws.Range("C13").Value2 = 14
$ws.Range("D13").Value2 = 265
$ws.Range("E13").Value2 = 14
$ws.Range("A14").Value2 = "Asturias"
$ws.Range("B14").Value2 = 292
$ws.Range("C14").Value2 = 4
$ws.Range("D14").Value2 = 287
$ws.Range("E14").Value2 = 1
$ws.Range("A15").Value2 = "Albacete"
$ws.Range("B15").Value2 = 259
$ws.Range("C15").Value2 = 8
$ws.Range("D15").Value2 = 231
$ws.Range("E15").Value2 = 20
$ws.Range("A16").Value2 = "Zaragoza"
$ws.Range("B16").Value2 = 224
$ws.Range("C16").Value2 = 0
$ws.Range("D16").Value2 = 210
$ws.Range("E16").Value2 = 14
$ws.Range("A17").Value2 = "A Coruña"
$ws.Range("B17").Value2 = 222
$ws.Range("C17").Value2 = 4
$ws.Range("D17").Value2 = 219
$ws.Range("E17").Value2 = 3
$ws.Range("A18").Value2 = "Ciudad Real"
$ws.Range("B18").Value2 = 216
$ws.Range("C18").Value2 = 6
$ws.Range("D18").Value2 = 189
$ws.Range("E18").Value2 = 21
$ws.Range("A19").Value2 = "Murcia"
$ws.Range("B19").Value2 = 215
$ws.Range("C19").Value2 = 1
$ws.Range("D19").Value2 = 214
$ws.Range("E19").Value2 = 0
$ws.Range("A20").Value2 = "Guadalajara"
$ws.Range("B20").Value2 = 205
$ws.Range("C20").Value2 = 2
$ws.Range("D20").Value2 = 200
$ws.Range("E20").Value2 = 3
$ws.Range("A21").Value2 = "Burgos"
$ws.Range("B21").Value2 = 187
$ws.Range("C21").Value2 = 14
$ws.Range("D21").Value2 = 165
$ws.Range("E21").Value2 = 8
$ws.Range("A22").Value2 = "Gipuzkoa/Guipuzcoa"
$ws.Range("B22").Value2 = 176
$ws.Range("C22").Value2 = 18
$ws.Range("D22").Value2 = 171
$ws.Range("E22").Value2 = 5
$ws.Range("A23").Value2 = "Granada"
$ws.Range("B23").Value2 = 176
$ws.Range("C23").Value2 = 0
$ws.Range("D23").Value2 = 169
$ws.Range("E23").Value2 = 7
$ws.Range("A24").Value2 = "Aragon"
$ws.Range("B24").Value2 = 174
$ws.Range("C24").Value2 = 0
$ws.Range("D24").Value2 = 163
$ws.Range("E24").Value2 = 11
$ws.Range("A25").Value2 = "Illes Balears"
$ws.Range("B25").Value2 = 169
$ws.Range("C25").Value2 = 6
$ws.Range("D25").Value2 = 161
$ws.Range("E25").Value2 = 2
$ws.Range("A26").Value2 = "Illes Balears*"
$ws.Range("B26").Value2 = 169
$ws.Range("C26").Value2 = 6
$ws.Range("D26").Value2 = 161
$ws.Range("E26").Value2 = 2
$ws.Range("A27").Value2 = "Caceres"
$ws.Range("B27").Value2 = 164
$ws.Range("C27").Value2 = 2
$ws.Range("D27").Value2 = 153
$ws.Range("E27").Value2 = 9
$ws.Range("A28").Value2 = "Salamanca"
$ws.Range("B28").Value2 = 149
$ws.Range("C28").Value2 = 8
$ws.Range("D28").Value2 = 129
$ws.Range("E28").Value2 = 12
$ws.Range("A29").Value2 = "Pontevedra"
$ws.Range("B29").Value2 = 145
$ws.Range("C29").Value2 = 4
$ws.Range("D29").Value2 = 143
$ws.Range("E29").Value2 = 2
$ws.Range("A30").Value2 = "Tenerife"
$ws.Range("B30").Value2 = 143
$ws.Range("C30").Value2 = 6
$ws.Range("D30").Value2 = 135
$ws.Range("E30").Value2 = 2
$ws.Range("A31").Value2 = "Leon"
$ws.Range("B31").Value2 = 134
$ws.Range("C31").Value2 = 3
$ws.Range("D31").Value2 = 126
$ws.Range("E31").Value2 = 5
$ws.Range("A32").Value2 = "Sevilla"
$ws.Range("B32").Value2 = 133
$ws.Range("C32").Value2 = 1
$ws.Range("D32").Value2 = 131
$ws.Range("E32").Value2 = 1
$ws.Range("A33").Value2 = "Segovia"
$ws.Range("B33").Value2 = 121
$ws.Range("C33").Value2 = 3
$ws.Range("D33").Value2 = 111
$ws.Range("E33").Value2 = 7
$ws.Range("A34").Value2 = "Valladolid"
$ws.Range("B34").Value2 = 115
$ws.Range("C34").Value2 = 1
$ws.Range("D34").Value2 = 111
$ws.Range("E34").Value2 = 3
$ws.Range("A35").Value2 = "Cordoba"
$ws.Range("B35").Value2 = 101
$ws.Range("C35").Value2 = 0
$ws.Range("D35").Value2 = 101
$ws.Range("E35").Value2 = 0
$ws.Range("A36").Value2 = "Jaen"
$ws.Range("B36").Value2 = 87
$ws.Range("C36").Value2 = 0
$ws.Range("D36").Value2 = 85
$ws.Range("E36").Value2 = 2
$ws.Range("A37").Value2 = "Cadiz"
$ws.Range("B37").Value2 = 84
$ws.Range("C37").Value2 = 0
$ws.Range("D37").Value2 = 84
$ws.Range("E37").Value2 = 0
$ws.Range("A38").Value2 = "Cantabria"
$ws.Range("B38").Value2 = 83
$ws.Range("C38").Value2 = 10
$ws.Range("D38").Value2 = 72
$ws.Range("E38").Value2 = 1
$ws.Range("A39").Value2 = "Badajoz"
$ws.Range("B39").Value2 = 77
$ws.Range("C39").Value2 = 4
$ws.Range("D39").Value2 = 73
$ws.Range("E39").Value2 = 0
$ws.Range("A40").Value2 = "Cuenca"
$ws.Range("B40").Value2 = 72
$ws.Range("C40").Value2 = 4
$ws.Range("D40").Value2 = 64
$ws.Range("E40").Value2 = 4
$ws.Range("A41").Value2 = "Castello/Castellon"
$ws.Range("B41").Value2 = 59
$ws.Range("C41").Value2 = 1
$ws.Range("D41").Value2 = 57
$ws.Range("E41").Value2 = 1
$ws.Range("A42").Value2 = "Avila"
$ws.Range("B42").Value2 = 59
$ws.Range("C42").Value2 = 2
$ws.Range("D42").Value2 = 55
$ws.Range("E42").Value2 = 2
$ws.Range("A43").Value2 = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Range("B43").Value2 = 58
$ws.Range("C43").Value2 = 0
$ws.Range("D43").Value2 = 58
$ws.Range("E43").Value2 = 3
$ws.Range("A44").Value2 = "Soria"
$ws.Range("B44").Value2 = 58
$ws.Range("C44").Value2 = 4
$ws.Range("D44").Value2 = 49
$ws.Range("E44").Value2 = 5
$ws.Range("A45").Value2 = "Gran Canaria"
$ws.Range("B45").Value2 = 55
$ws.Range("C45").Value2 = 0
$ws.Range("D45").Value2 = 54
$ws.Range("E45").Value2 = 1
$ws.Range("A46").Value2 = "Ourense"
$ws.Range("B46").Value2 = 46
$ws.Range("C46").Value2 = 4
$ws.Range("D46").Value2 = 46
$ws.Range("E46").Value2 = 0
$ws.Range("A47").Value2 = "Almeria"
$ws.Range("B47").Value2 = 37
$ws.Range("C47").Value2 = 0
$ws.Range("D47").Value2 = 37
$ws.Range("E47").Value2 = 0
$ws.Range("A48").Value2 = "Lugo"
$ws.Range("B48").Value2 = 36
$ws.Range("C48").Value2 = 4
$ws.Range("D48").Value2 = 36
$ws.Range("E48").Value2 = 0
$ws.Range("A49").Value2 = "Zamora"
$ws.Range("B49").Value2 = 31
$ws.Range("C49").Value2 = 1
$ws.Range("D49").Value2 = 29
$ws.Range("E49").Value2 = 1
$ws.Range("A50").Value2 = "Teruel"
$ws.Range("B50").Value2 = 27
$ws.Range("C50").Value2 = 0
$ws.Range("D50").Value2 = 26
$ws.Range("E50").Value2 = 1
$ws.Range("A51").Value2 = "Huesca"
$ws.Range("B51").Value2 = 24
$ws.Range("C51").Value2 = 0
$ws.Range("D51").Value2 = 24
$ws.Range("E51").Value2 = 0
$ws.Range("A52").Value2 = "Melilla"
$ws.Range("B52").Value2 = 23
$ws.Range("C52").Value2 = 0
$ws.Range("D52").Value2 = 23
$ws.Range("E52").Value2 = 0
$ws.Range("A53").Value2 = "Huelva"
$ws.Range("B53").Value2 = 23
$ws.Range("C53").Value2 = 0
$ws.Range("D53").Value2 = 23
$ws.Range("E53").Value2 = 0
$ws.Range("A54").Value2 = "Palencia"
$ws.Range("B54").Value2 = 14
$ws.Range("C54").Value2 = 1
$ws.Range("D54").Value2 = 13
$ws.Range("E54").Value2 = 0
$ws.Range("A55").Value2 = "Fuerteventura"
$ws.Range("B55").Value2 = 11
$ws.Range("C55").Value2 = 0
$ws.Range("D55").Value2 = 11
$ws.Range("E55").Value2 = 0
$ws.Range("A56").Value2 = "Arroyo de la Luz"
$ws.Range("B56").Value2 = 7
$ws.Range("C56").Value2 = 0
$ws.Range("D56").Value2 = 7
$ws.Range("E56").Value2 = 0
$ws.Range("A57").Value2 = "La Palma"
$ws.Range("B57").Value2 = 5
$ws.Range("C57").Value2 = 0
$ws.Range("D57").Value2 = 5
$ws.Range("E57").Value2 = 0
$ws.Range("A58").Value2 = "Ceuta"
$ws.Range("B58").Value2 = 5
$ws.Range("C58").Value2 = 0
$ws.Range("D58").Value2 = 5
$ws.Range("E58").Value2 = 0
$ws.Range("A59").Value2 = "Lanzarote"
$ws.Range("B59").Value2 = 3
$ws.Range("C59").Value2 = 0
$ws.Range("D59").Value2 = 3
$ws.Range("E59").Value2 = 0
$ws.Range("A60").Value2 = "La Gomera"
$ws.Range("B60").Value2 = 3
$ws.Range("C60").Value2 = 2
$ws.Range("D60").Value2 = 1
$ws.Range("E60").Value2 = 0

# The source had one extra trailing row (61, "La Gomera" duplicate slot) that
# is no longer needed now that the duplicate Cataluna* row was folded in above;
# delete it so the used range shrinks back to A1:E60.
$ws.Rows(61).Delete()
